$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.181.49"
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").Value = "1.849.53"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'0.7055"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.17%  "
$ws.Range("D6").Value = "'238.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.3054"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").Value = "'0.07425"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.65%  "
$ws.Range("D10").Value = "'23.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.37%  "
$ws.Range("D11").Value = "'0.08136"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.32%  "
$ws.Range("D12").Value = "'0.7280"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.51%  "
$ws.Range("D13").Value = "1.843.79"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").Value = "'5.236"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "'88.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").Value = "29.171.53"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("D17").Value = "'5.774"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.86%  "
$ws.Range("D18").Value = "'238.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.31%  "
$ws.Range("D19").Value = "'13.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "'0.000007647"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").Value = "'0.9998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "2.095.05"
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'7.607"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.60%  "
$ws.Range("D25").Value = "'9.007"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("D26").Value = "'161.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "'0.1457"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.94%  "
$ws.Range("D28").Value = "'18.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").Value = "'1.971"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.54%  "
$ws.Range("E30").Value = "  -5.20%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'1.493"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("D33").Value = "'3.991"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.45%  "
$ws.Range("D34").Value = "'0.05185"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("D35").Value = "'1.188"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.73%  "
$ws.Range("D36").Value = "'1.032"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("D37").Value = "'0.7046"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.30%  "
$ws.Range("D38").Value = "'2.662"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("D39").Value = "'0.01870"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.20%  "
$ws.Range("D40").Value = "'2.683"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").Value = "'0.9392"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.39%  "
$ws.Range("D42").Value = "'6.014"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").Value = "1.075.49"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").Value = "'0.4294"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.54%  "
$ws.Range("D45").Value = "'70.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").Value = "'0.9997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'102.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "'1.745"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.69%  "
$ws.Range("D49").Value = "1.986.13"
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("D50").Value = "'7.067"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.14%  "
$ws.Range("D51").Value = "'9.107"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.90%  "
